# Update marksheet correction/total values on the "quiz" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# B11: marking score 3 -> 5
$ws.Range("B11").Value = 5

# B12: total score 36 -> 60
$ws.Range("B12").Value = 60

# E12: correct/total text "33/84" -> "60/140"
$ws.Range("E12").Value = "60/140"
